$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 34835
$ws.Cells.Item(2, 4).Value = 50453444
$ws.Cells.Item(3, 3).Value = 85403
$ws.Cells.Item(3, 4).Value = 125351861
$ws.Cells.Item(4, 3).Value = 29297
$ws.Cells.Item(4, 4).Value = 43430308
$ws.Cells.Item(5, 3).Value = 8064
$ws.Cells.Item(5, 4).Value = 11990518
$ws.Cells.Item(6, 3).Value = 1731
$ws.Cells.Item(6, 4).Value = 2573596
$ws.Cells.Item(11, 3).Value = 38266
$ws.Cells.Item(11, 4).Value = 52064898
$ws.Cells.Item(12, 3).Value = 8969
$ws.Cells.Item(12, 4).Value = 12987230
$ws.Cells.Item(13, 3).Value = 24594
$ws.Cells.Item(13, 4).Value = 36096298
$ws.Cells.Item(14, 3).Value = 7845
$ws.Cells.Item(14, 4).Value = 11653450
$ws.Cells.Item(16, 3).Value = 362
$ws.Cells.Item(16, 4).Value = 532123
$ws.Cells.Item(19, 3).Value = 9498
$ws.Cells.Item(19, 4).Value = 12628286
$ws.Cells.Item(20, 3).Value = 12529
$ws.Cells.Item(20, 4).Value = 18105522
$ws.Cells.Item(21, 3).Value = 30017
$ws.Cells.Item(21, 4).Value = 44105912
$ws.Cells.Item(22, 3).Value = 9742
$ws.Cells.Item(22, 4).Value = 14489310
$ws.Cells.Item(23, 3).Value = 2469
$ws.Cells.Item(23, 4).Value = 3673763
$ws.Cells.Item(26, 3).Value = 10947
$ws.Cells.Item(26, 4).Value = 14688423
$ws.Cells.Item(28, 3).Value = 21254
$ws.Cells.Item(28, 4).Value = 31227892
$ws.Cells.Item(29, 3).Value = 7390
$ws.Cells.Item(29, 4).Value = 11000651
$ws.Cells.Item(30, 3).Value = 1821
$ws.Cells.Item(30, 4).Value = 2719986
$ws.Cells.Item(31, 3).Value = 298
$ws.Cells.Item(31, 4).Value = 444415
$ws.Cells.Item(33, 3).Value = 7739
$ws.Cells.Item(33, 4).Value = 10262641
$ws.Cells.Item(34, 3).Value = 2778
$ws.Cells.Item(34, 4).Value = 4007836
$ws.Cells.Item(35, 3).Value = 7000
$ws.Cells.Item(35, 4).Value = 10230347
$ws.Cells.Item(36, 3).Value = 2838
$ws.Cells.Item(36, 4).Value = 4202023
$ws.Cells.Item(40, 3).Value = 2159
$ws.Cells.Item(40, 4).Value = 2911166
$ws.Cells.Item(41, 3).Value = 16058
$ws.Cells.Item(41, 4).Value = 23241051
$ws.Cells.Item(42, 3).Value = 48231
$ws.Cells.Item(42, 4).Value = 70766433
$ws.Cells.Item(43, 3).Value = 18067
$ws.Cells.Item(43, 4).Value = 26843123
$ws.Cells.Item(44, 3).Value = 5243
$ws.Cells.Item(44, 4).Value = 7815108
$ws.Cells.Item(49, 3).Value = 15507
$ws.Cells.Item(49, 4).Value = 20715479
$ws.Cells.Item(50, 3).Value = 1682
$ws.Cells.Item(50, 4).Value = 2440833
$ws.Cells.Item(51, 3).Value = 5976
$ws.Cells.Item(51, 4).Value = 8799326
$ws.Cells.Item(52, 3).Value = 2090
$ws.Cells.Item(52, 4).Value = 3121750
$ws.Cells.Item(53, 3).Value = 683
$ws.Cells.Item(53, 4).Value = 1019805
$ws.Cells.Item(54, 3).Value = 153
$ws.Cells.Item(54, 4).Value = 227111
$ws.Cells.Item(56, 3).Value = 5483
$ws.Cells.Item(56, 4).Value = 7572623
$ws.Cells.Item(57, 3).Value = 690
$ws.Cells.Item(57, 4).Value = 1011340
$ws.Cells.Item(58, 3).Value = 1743
$ws.Cells.Item(58, 4).Value = 2582914
$ws.Cells.Item(61, 3).Value = 51
$ws.Cells.Item(61, 4).Value = 76500
$ws.Cells.Item(63, 3).Value = 1008
$ws.Cells.Item(63, 4).Value = 1429765
$ws.Cells.Item(64, 3).Value = 14297
$ws.Cells.Item(64, 4).Value = 20669944
$ws.Cells.Item(65, 3).Value = 42359
$ws.Cells.Item(65, 4).Value = 62036210
$ws.Cells.Item(66, 3).Value = 14950
$ws.Cells.Item(66, 4).Value = 22230375
$ws.Cells.Item(67, 3).Value = 4298
$ws.Cells.Item(67, 4).Value = 6403793
$ws.Cells.Item(68, 3).Value = 829
$ws.Cells.Item(68, 4).Value = 1234273
$ws.Cells.Item(71, 3).Value = 14184
$ws.Cells.Item(71, 4).Value = 18782082
$ws.Cells.Item(72, 3).Value = 45574
$ws.Cells.Item(72, 4).Value = 66353129
$ws.Cells.Item(73, 3).Value = 132349
$ws.Cells.Item(73, 4).Value = 195116962
$ws.Cells.Item(74, 3).Value = 58244
$ws.Cells.Item(74, 4).Value = 86820399
$ws.Cells.Item(75, 3).Value = 18482
$ws.Cells.Item(75, 4).Value = 27620072
$ws.Cells.Item(76, 3).Value = 4121
$ws.Cells.Item(76, 4).Value = 6158120
$ws.Cells.Item(83, 3).Value = 45042
$ws.Cells.Item(83, 4).Value = 61630269
$ws.Cells.Item(84, 3).Value = 4082
$ws.Cells.Item(84, 4).Value = 5920669
$ws.Cells.Item(85, 3).Value = 10585
$ws.Cells.Item(85, 4).Value = 15559960
$ws.Cells.Item(86, 3).Value = 3624
$ws.Cells.Item(86, 4).Value = 5402081
$ws.Cells.Item(92, 3).Value = 1385
$ws.Cells.Item(92, 4).Value = 2001990
$ws.Cells.Item(93, 3).Value = 4588
$ws.Cells.Item(93, 4).Value = 6761465
$ws.Cells.Item(94, 3).Value = 1778
$ws.Cells.Item(94, 4).Value = 2651616
$ws.Cells.Item(99, 3).Value = 3041
$ws.Cells.Item(99, 4).Value = 4037232
$ws.Cells.Item(100, 3).Value = 517
$ws.Cells.Item(100, 4).Value = 770964
$ws.Cells.Item(101, 3).Value = 266
$ws.Cells.Item(101, 4).Value = 397265
$ws.Cells.Item(105, 3).Value = 10043
$ws.Cells.Item(105, 4).Value = 14593858
$ws.Cells.Item(106, 3).Value = 27756
$ws.Cells.Item(106, 4).Value = 40808961
$ws.Cells.Item(108, 3).Value = 2530
$ws.Cells.Item(108, 4).Value = 3772410
$ws.Cells.Item(109, 3).Value = 432
$ws.Cells.Item(109, 4).Value = 645482
$ws.Cells.Item(112, 3).Value = 9150
$ws.Cells.Item(112, 4).Value = 12129717
$ws.Cells.Item(113, 3).Value = 28202
$ws.Cells.Item(113, 4).Value = 40715452
$ws.Cells.Item(114, 3).Value = 62522
$ws.Cells.Item(114, 4).Value = 91579750
$ws.Cells.Item(115, 3).Value = 20271
$ws.Cells.Item(115, 4).Value = 30140526
$ws.Cells.Item(116, 3).Value = 5671
$ws.Cells.Item(116, 4).Value = 8453861
$ws.Cells.Item(117, 3).Value = 1014
$ws.Cells.Item(117, 4).Value = 1516993
$ws.Cells.Item(118, 3).Value = 59
$ws.Cells.Item(118, 4).Value = 85920
$ws.Cells.Item(121, 3).Value = 24047
$ws.Cells.Item(121, 4).Value = 32215789
$ws.Cells.Item(122, 3).Value = 32874
$ws.Cells.Item(122, 4).Value = 47502784
$ws.Cells.Item(123, 3).Value = 71695
$ws.Cells.Item(123, 4).Value = 104952007
$ws.Cells.Item(124, 3).Value = 22444
$ws.Cells.Item(124, 4).Value = 33325154
$ws.Cells.Item(125, 3).Value = 5930
$ws.Cells.Item(125, 4).Value = 8818054
$ws.Cells.Item(126, 3).Value = 1081
$ws.Cells.Item(126, 4).Value = 1610551
$ws.Cells.Item(130, 3).Value = 29157
$ws.Cells.Item(130, 4).Value = 38850427
$ws.Cells.Item(131, 3).Value = 12319
$ws.Cells.Item(131, 4).Value = 17844177
$ws.Cells.Item(132, 3).Value = 30669
$ws.Cells.Item(132, 4).Value = 45079415
$ws.Cells.Item(133, 3).Value = 10932
$ws.Cells.Item(133, 4).Value = 16244690
$ws.Cells.Item(138, 3).Value = 10131
$ws.Cells.Item(138, 4).Value = 13572831
$ws.Cells.Item(139, 3).Value = 32128
$ws.Cells.Item(139, 4).Value = 46433342
$ws.Cells.Item(140, 3).Value = 75922
$ws.Cells.Item(140, 4).Value = 111305711
$ws.Cells.Item(141, 3).Value = 22915
$ws.Cells.Item(141, 4).Value = 34076653
$ws.Cells.Item(142, 3).Value = 5905
$ws.Cells.Item(142, 4).Value = 8814848
$ws.Cells.Item(143, 3).Value = 1277
$ws.Cells.Item(143, 4).Value = 1902186
$ws.Cells.Item(146, 3).Value = 27034
$ws.Cells.Item(146, 4).Value = 36639603
